$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the cell values and the shared
# string text can be updated, then re-protect at the end.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A41).
$ws.Cells.Item(41, 1).Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-29 for illustrative purposes only and are subject to change."

# Updated Weight (D) / Percent Change (E) figures for rows 2-38.
$data = @(
    @{ Row = 2; D = 0.0303196701138163; E = -0.009291521486643362 },
    @{ Row = 3; D = 0.02992147557198382; E = 0.001183363079071853 },
    @{ Row = 4; D = 0.03088774551460356; E = -0.001309212092358902 },
    @{ Row = 5; D = 0.06691719256208196; E = 0.00370391788347546 },
    @{ Row = 6; D = 0.01505012839742369; E = 0.0230638691761802 },
    @{ Row = 7; D = 0.01585464390030973; E = 0.04473896170462055 },
    @{ Row = 8; D = 0.02748006705456774; E = -0.00844915719656969 },
    @{ Row = 9; D = 0.03445214777390289; E = 0.005279119397955778 },
    @{ Row = 10; D = 0.02940138474183527; E = 0.004896154149885534 },
    @{ Row = 11; D = 0.03206143263353254; E = -0.01072997634335926 },
    @{ Row = 12; D = 0.01345212610345164; E = -0.02711254944264652 },
    @{ Row = 13; D = 0.0141020461549654; E = 0.01325393776411832 },
    @{ Row = 14; D = 0.01631514098950376; E = -0.02805910675742984 },
    @{ Row = 15; D = 0.008631031157465236; E = 0.02555595408895273 },
    @{ Row = 16; D = 0.007493138980343352; E = 0.03276783639321401 },
    @{ Row = 17; D = 0.03208310308478873; E = -0.01063829787234039 },
    @{ Row = 18; D = 0.02947974664146702; E = -0.0757346039997111 },
    @{ Row = 19; D = 0.03162918452544033; E = 0.009359515507432636 },
    @{ Row = 20; D = 0.02970980169989211; E = 0.07297297297297289 },
    @{ Row = 21; D = 0.04604796754385557; E = 0.02100079414767797 },
    @{ Row = 22; D = 0.03367723565532219; E = 0.01396110424865693 },
    @{ Row = 23; D = 0.03092973201391242; E = 0.01760970879859869 },
    @{ Row = 24; D = 0.02941531574621424; E = 0.01835188254795184 },
    @{ Row = 25; D = 0.01523181191286621; E = -0.02938150222933578 },
    @{ Row = 26; D = 0.01512326617041333; E = -0.009518691947493774 },
    @{ Row = 27; D = 0.0304353748446306; E = 0.008671328671328693 },
    @{ Row = 28; D = 0.02909103292205912; E = 0.01239757369373207 },
    @{ Row = 29; D = 0.02955230395594088; E = -0.008053111250785783 },
    @{ Row = 30; D = 0.02781963528630536; E = 0.01950883635529022 },
    @{ Row = 31; D = 0.03547007852442801; E = 0.003142029554715542 },
    @{ Row = 32; D = 0.03018635814135635; E = 0.01487696538108985 },
    @{ Row = 33; D = 0.03009832193312808; E = -0.004679926458298511 },
    @{ Row = 34; D = 0.03053908343278522; E = 0.01763856154489485 },
    @{ Row = 35; D = 0.03003872819217356; E = 0.0002318840579709214 },
    @{ Row = 36; D = 0.02945439995294416; E = 0.0194442619720161 },
    @{ Row = 37; D = 0.0316481461702895; E = 0.01428152205810429 },
    @{ Row = 38; D = 0.9999999999999999; E = 0.005337316052696162 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

$ws.Protect()
